$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 18: item 5471
$ws.Range("H18").Value = 311
$ws.Range("I18").Value = 311
$ws.Range("K18").Value = 311
$ws.Range("M18").Value = -27

# ALC row 19: item 7015
$ws.Range("H19").Value = 1149.7142
$ws.Range("I19").Value = 1909
$ws.Range("J19").Value = 390.42856
$ws.Range("K19").Value = 1909
$ws.Range("L19").Value = 390.42856
$ws.Range("M19").Value = -1734
$ws.Range("N19").Value = -740.4285600000001

# ALC row 32: item 5484
$ws.Range("H32").Value = 7860
$ws.Range("J32").Value = 8700
$ws.Range("L32").Value = 8700
$ws.Range("N32").Value = -9352

# ALC row 33: item 5512
$ws.Range("H33").Value = 142.14285
$ws.Range("I33").Value = 149.18182
$ws.Range("J33").Value = 116.333336
$ws.Range("K33").Value = 149.18182
$ws.Range("L33").Value = 116.333336
$ws.Range("M33").Value = 79.81818000000001
$ws.Range("N33").Value = -574.333336

# ALC row 43: item 5472
$ws.Range("H43").Value = 1616.6666
$ws.Range("J43").Value = 1925
$ws.Range("L43").Value = 1925
$ws.Range("N43").Value = -2063

# ALC row 74: item 5507
$ws.Range("H74").Value = 5359.6
$ws.Range("I74").Value = 4599.6665
$ws.Range("J74").Value = 6499.5
$ws.Range("K74").Value = 4599.6665
$ws.Range("L74").Value = 6499.5
$ws.Range("M74").Value = -3663.6665
$ws.Range("N74").Value = -8371.5

# ALC row 77: item 5507
$ws.Range("H77").Value = 5359.6
$ws.Range("I77").Value = 4599.6665
$ws.Range("J77").Value = 6499.5
$ws.Range("K77").Value = 22998.3325
$ws.Range("L77").Value = 32497.5
$ws.Range("M77").Value = -18318.3325
$ws.Range("N77").Value = -41857.5

# ALC row 113: item 27775
$ws.Range("H113").Value = 2620.5715
$ws.Range("I113").Value = 2265.6667
$ws.Range("K113").Value = 2265.6667
$ws.Range("M113").Value = 988.3332999999998

# ALC row 137: item 44013
$ws.Range("H137").Value = 4400.2
$ws.Range("I137").Value = 5249.5
$ws.Range("K137").Value = 15748.5
$ws.Range("M137").Value = -13198.5

# ALC row 138: item 44169
$ws.Range("H138").Value = 5886008.5
$ws.Range("J138").Value = 3866.25
$ws.Range("L138").Value = 11598.75
$ws.Range("N138").Value = -21878.75

$ws = $wb.Worksheets.Item("ARM")
# ARM row 74: item 44000
$ws.Range("H74").Value = 7968.6665
$ws.Range("I74").Value = 2956
$ws.Range("J74").Value = 10475
$ws.Range("K74").Value = 2956
$ws.Range("L74").Value = 10475
$ws.Range("M74").Value = -2082
$ws.Range("N74").Value = -12223

# ARM row 77: item 44000
$ws.Range("H77").Value = 7968.6665
$ws.Range("I77").Value = 2956
$ws.Range("J77").Value = 10475
$ws.Range("K77").Value = 14780
$ws.Range("L77").Value = 52375
$ws.Range("M77").Value = -10412
$ws.Range("N77").Value = -61111

# ARM row 97: item 19941
$ws.Range("H97").Value = 869.3333
$ws.Range("I97").Value = 869.3333
$ws.Range("K97").Value = 869.3333
$ws.Range("M97").Value = -373.3333

$ws = $wb.Worksheets.Item("BSM")
# BSM row 37: item 2485
$ws.Range("H37").Value = 3075
$ws.Range("I37").Value = 650
$ws.Range("J37").Value = 5500
$ws.Range("K37").Value = 650
$ws.Range("L37").Value = 5500
$ws.Range("M37").Value = -513
$ws.Range("N37").Value = -5774

# BSM row 94: item 19939
$ws.Range("H94").Value = 1336.5
$ws.Range("I94").Value = 1336.5
$ws.Range("K94").Value = 1336.5
$ws.Range("M94").Value = -885.5

# BSM row 99: item 19943
$ws.Range("H99").Value = 1129.3334
$ws.Range("I99").Value = 1129.3334
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1129.3334
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 368.6666
$ws.Range("N99").ClearContents()

# BSM row 102: item 19565
$ws.Range("H102").Value = 25733.334
$ws.Range("I102").Value = 15880
$ws.Range("K102").Value = 15880
$ws.Range("M102").Value = -12635

# BSM row 105: item 19947
$ws.Range("H105").Value = 2099.3333
$ws.Range("I105").Value = 2049.25
$ws.Range("K105").Value = 2049.25
$ws.Range("M105").Value = -302.25

$ws = $wb.Worksheets.Item("CRP")
# CRP row 31: item 44023
$ws.Range("H31").Value = 6148.0625
$ws.Range("I31").Value = 2107.7778
$ws.Range("K31").Value = 2107.7778
$ws.Range("M31").Value = -1812.7778

# CRP row 34: item 44023
$ws.Range("H34").Value = 6148.0625
$ws.Range("I34").Value = 2107.7778
$ws.Range("K34").Value = 2107.7778
$ws.Range("M34").Value = -1905.7778

# CRP row 58: item 44021
$ws.Range("H58").Value = 6647.7334
$ws.Range("I58").Value = 3772.2
$ws.Range("K58").Value = 3772.2
$ws.Range("M58").Value = -3569.2

# CRP row 105: item 19928
$ws.Range("H105").Value = 3999
$ws.Range("I105").Value = 3999
$ws.Range("K105").Value = 3999
$ws.Range("M105").Value = -2252

# CRP row 107: item 27689
$ws.Range("H107").Value = 662.55554
$ws.Range("I107").Value = 620.375
$ws.Range("K107").Value = 620.375
$ws.Range("M107").Value = 1299.625

# CRP row 136: item 44021
$ws.Range("H136").Value = 6647.7334
$ws.Range("I136").Value = 3772.2
$ws.Range("K136").Value = 11316.6
$ws.Range("M136").Value = -8766.599999999999

$ws = $wb.Worksheets.Item("CUL")
# CUL row 64: item 12861
$ws.Range("H64").Value = 7998
$ws.Range("I64").Value = 7998
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 23994
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = -23724
$ws.Range("N64").ClearContents()

# CUL row 67: item 12861
$ws.Range("H67").Value = 7998
$ws.Range("I67").Value = 7998
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 23994
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = -23058
$ws.Range("N67").ClearContents()

# CUL row 86: item 12892
$ws.Range("H86").Value = 93.40000000000001
$ws.Range("I86").Value = 100
$ws.Range("J86").Value = 91.75
$ws.Range("K86").Value = 300
$ws.Range("L86").Value = 275.25
$ws.Range("M86").Value = 886
$ws.Range("N86").Value = -2647.25

# CUL row 89: item 12892
$ws.Range("H89").Value = 93.40000000000001
$ws.Range("I89").Value = 100
$ws.Range("J89").Value = 91.75
$ws.Range("K89").Value = 900
$ws.Range("L89").Value = 825.75
$ws.Range("M89").Value = 5028
$ws.Range("N89").Value = -12681.75

$ws = $wb.Worksheets.Item("GSM")
# GSM row 18: item 4309
$ws.Range("H18").Value = 10000
$ws.Range("J18").Value = 10000
$ws.Range("L18").Value = 10000
$ws.Range("N18").Value = -10586

# GSM row 21: item 4430
$ws.Range("H21").Value = 10000
$ws.Range("J21").Value = 10000
$ws.Range("L21").Value = 10000
$ws.Range("N21").Value = -10346

# GSM row 30: item 4430
$ws.Range("H30").Value = 10000
$ws.Range("J30").Value = 10000
$ws.Range("L30").Value = 10000
$ws.Range("N30").Value = -10210

# GSM row 41: item 2449
$ws.Range("H41").Value = 4811.8
$ws.Range("I41").Value = 3499.75
$ws.Range("K41").Value = 3499.75
$ws.Range("M41").Value = -3144.75

# GSM row 44: item 4143
$ws.Range("H44").Value = 30030
$ws.Range("J44").Value = 30030
$ws.Range("L44").Value = 30030
$ws.Range("N44").Value = -31222

# GSM row 80: item 12521
$ws.Range("H80").Value = 3969.75
$ws.Range("J80").Value = 3969.75
$ws.Range("L80").Value = 3969.75
$ws.Range("N80").Value = -5965.75

# GSM row 83: item 12521
$ws.Range("H83").Value = 3969.75
$ws.Range("J83").Value = 3969.75
$ws.Range("L83").Value = 19848.75
$ws.Range("N83").Value = -29832.75

# GSM row 97: item 19940
$ws.Range("H97").Value = 848
$ws.Range("I97").Value = 700
$ws.Range("J97").Value = 996
$ws.Range("K97").Value = 700
$ws.Range("L97").Value = 996
$ws.Range("M97").Value = -204
$ws.Range("N97").Value = -1988

# GSM row 132: item 44008
$ws.Range("H132").Value = 3248.2632
$ws.Range("I132").Value = 1558.7858
$ws.Range("J132").Value = 7978.8
$ws.Range("K132").Value = 4676.357400000001
$ws.Range("L132").Value = 23936.4
$ws.Range("M132").Value = -2146.357400000001
$ws.Range("N132").Value = -28996.4

$ws = $wb.Worksheets.Item("LTW")
# LTW row 13: item 3546
$ws.Range("H13").Value = 18666
$ws.Range("J13").Value = 17999
$ws.Range("L13").Value = 17999
$ws.Range("N13").Value = -18279

# LTW row 22: item 5277
$ws.Range("H22").Value = 3933
$ws.Range("I22").Value = 899.5
$ws.Range("K22").Value = 899.5
$ws.Range("M22").Value = -604.5

# LTW row 27: item 5277
$ws.Range("H27").Value = 3933
$ws.Range("I27").Value = 899.5
$ws.Range("K27").Value = 899.5
$ws.Range("M27").Value = -792.5

# LTW row 55: item 5284
$ws.Range("H55").Value = 390
$ws.Range("I55").Value = 425
$ws.Range("K55").Value = 425
$ws.Range("M55").Value = -252

# LTW row 100: item 19995
$ws.Range("H100").Value = 1850
$ws.Range("I100").Value = 1850
$ws.Range("K100").Value = 1850
$ws.Range("M100").Value = -1309

$ws = $wb.Worksheets.Item("WVR")
# WVR row 132: item 44029
$ws.Range("H132").Value = 6253.2856
$ws.Range("I132").Value = 943.25
$ws.Range("J132").Value = 13333.333
$ws.Range("K132").Value = 2829.75
$ws.Range("L132").Value = 39999.999
$ws.Range("M132").Value = -299.75
$ws.Range("N132").Value = -45059.999
